$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data rows (2-12) with placeholder column-name tokens,
# matching the new database-driven OrderManager template.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "sku"
    $ws.Cells.Item($r, 2).Value = "name"
    $ws.Cells.Item($r, 3).Value = "quantity"
    $ws.Cells.Item($r, 4).Value = "cost_per"
    $ws.Cells.Item($r, 5).Value = "total_cost"
}
